$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq $null -or $v -eq "") { continue }

    $parts = $v -split ', '
    $others = @()
    $systemCount = 0
    foreach ($p in $parts) {
        if ($p.CompareTo("System") -eq 0) {
            $systemCount = $systemCount + 1
        } else {
            $others += $p
        }
    }
    if ($systemCount -eq 0) { continue }

    $newParts = @()
    $newParts += $others
    for ($i = 0; $i -lt $systemCount; $i++) {
        $newParts += "System"
    }
    $newValue = $newParts -join ', '

    if ($newValue -ne $v) {
        $cell.Value = $newValue
    }
}

Write-Host "Done. LastRow=$lastRow"
